$wb = $excel.ActiveWorkbook

# ---- Sheet: Triple Dribble ----
$ws = $wb.Worksheets.Item("Triple Dribble")
$ws.Range("A49:N49").Copy()
$ws.Range("A50:N52").PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Cells.Item(50,7).PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Cells.Item(51,7).PasteSpecial(-4122)
$ws.Range("G5").Copy()
$ws.Cells.Item(52,7).PasteSpecial(-4122)

$ws.Cells.Item(50,1).Value = "LUMI"
$ws.Cells.Item(50,2).Value = "CROW"
$ws.Cells.Item(50,3).Value = "MEG"
$ws.Cells.Item(50,4).Value = "CORDELIUS"
$ws.Cells.Item(50,5).Value = "JAE-YONG"
$ws.Cells.Item(50,6).Value = "HANK"
$ws.Cells.Item(50,7).Value = "Equipo 2"
$ws.Cells.Item(50,8).Value = "NOVO|Marco"
$ws.Cells.Item(50,9).Value = "TTM|Maury"
$ws.Cells.Item(50,10).Value = "TTM|Maru"
$ws.Cells.Item(50,11).Value = "Enraged 💔"
$ws.Cells.Item(50,12).Value = "Drage🍥"
$ws.Cells.Item(50,13).Value = "SUP|Tomzy"
$ws.Cells.Item(50,14).Value = "20250724T201248.000Z"
$ws.Cells.Item(51,1).Value = "LUMI"
$ws.Cells.Item(51,2).Value = "CROW"
$ws.Cells.Item(51,3).Value = "MEG"
$ws.Cells.Item(51,4).Value = "CORDELIUS"
$ws.Cells.Item(51,5).Value = "JAE-YONG"
$ws.Cells.Item(51,6).Value = "HANK"
$ws.Cells.Item(51,7).Value = "Equipo 1"
$ws.Cells.Item(51,8).Value = "NOVO|Marco"
$ws.Cells.Item(51,9).Value = "TTM|Maury"
$ws.Cells.Item(51,10).Value = "TTM|Maru"
$ws.Cells.Item(51,11).Value = "Enraged 💔"
$ws.Cells.Item(51,12).Value = "Drage🍥"
$ws.Cells.Item(51,13).Value = "SUP|Tomzy"
$ws.Cells.Item(51,14).Value = "20250724T201056.000Z"
$ws.Cells.Item(52,1).Value = "LUMI"
$ws.Cells.Item(52,2).Value = "CROW"
$ws.Cells.Item(52,3).Value = "MEG"
$ws.Cells.Item(52,4).Value = "CORDELIUS"
$ws.Cells.Item(52,5).Value = "JAE-YONG"
$ws.Cells.Item(52,6).Value = "HANK"
$ws.Cells.Item(52,7).Value = "Equipo 2"
$ws.Cells.Item(52,8).Value = "NOVO|Marco"
$ws.Cells.Item(52,9).Value = "TTM|Maury"
$ws.Cells.Item(52,10).Value = "TTM|Maru"
$ws.Cells.Item(52,11).Value = "Enraged 💔"
$ws.Cells.Item(52,12).Value = "Drage🍥"
$ws.Cells.Item(52,13).Value = "SUP|Tomzy"
$ws.Cells.Item(52,14).Value = "20250724T200846.000Z"

# ---- Sheet: Dueling Beetles ----
$ws = $wb.Worksheets.Item("Dueling Beetles")
$ws.Range("A21:N21").Copy()
$ws.Range("A22:N24").PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Cells.Item(22,7).PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Cells.Item(23,7).PasteSpecial(-4122)
$ws.Range("G6").Copy()
$ws.Cells.Item(24,7).PasteSpecial(-4122)

$ws.Cells.Item(22,1).Value = "CROW"
$ws.Cells.Item(22,2).Value = "BARLEY"
$ws.Cells.Item(22,3).Value = "KIT"
$ws.Cells.Item(22,4).Value = "GUS"
$ws.Cells.Item(22,5).Value = "MOE"
$ws.Cells.Item(22,6).Value = "FRANK"
$ws.Cells.Item(22,7).Value = "Equipo 1"
$ws.Cells.Item(22,8).Value = "NHG|Xemp"
$ws.Cells.Item(22,9).Value = "KCP|Tyrant"
$ws.Cells.Item(22,10).Value = "KCP|Fade"
$ws.Cells.Item(22,11).Value = "TRB|Zeus 解開"
$ws.Cells.Item(22,12).Value = "TRB|Lxffy"
$ws.Cells.Item(22,13).Value = "TRB|R B M"
$ws.Cells.Item(22,14).Value = "20250724T201924.000Z"
$ws.Cells.Item(23,1).Value = "CROW"
$ws.Cells.Item(23,2).Value = "BARLEY"
$ws.Cells.Item(23,3).Value = "KIT"
$ws.Cells.Item(23,4).Value = "GUS"
$ws.Cells.Item(23,5).Value = "MOE"
$ws.Cells.Item(23,6).Value = "FRANK"
$ws.Cells.Item(23,7).Value = "Equipo 2"
$ws.Cells.Item(23,8).Value = "NHG|Xemp"
$ws.Cells.Item(23,9).Value = "KCP|Tyrant"
$ws.Cells.Item(23,10).Value = "KCP|Fade"
$ws.Cells.Item(23,11).Value = "TRB|Zeus 解開"
$ws.Cells.Item(23,12).Value = "TRB|Lxffy"
$ws.Cells.Item(23,13).Value = "TRB|R B M"
$ws.Cells.Item(23,14).Value = "20250724T201650.000Z"
$ws.Cells.Item(24,1).Value = "CROW"
$ws.Cells.Item(24,2).Value = "BARLEY"
$ws.Cells.Item(24,3).Value = "KIT"
$ws.Cells.Item(24,4).Value = "GUS"
$ws.Cells.Item(24,5).Value = "MOE"
$ws.Cells.Item(24,6).Value = "FRANK"
$ws.Cells.Item(24,7).Value = "Equipo 1"
$ws.Cells.Item(24,8).Value = "NHG|Xemp"
$ws.Cells.Item(24,9).Value = "KCP|Tyrant"
$ws.Cells.Item(24,10).Value = "KCP|Fade"
$ws.Cells.Item(24,11).Value = "TRB|Zeus 解開"
$ws.Cells.Item(24,12).Value = "TRB|Lxffy"
$ws.Cells.Item(24,13).Value = "TRB|R B M"
$ws.Cells.Item(24,14).Value = "20250724T201500.000Z"

# ---- Sheet: Hard Rock Mine ----
$ws = $wb.Worksheets.Item("Hard Rock Mine")
$ws.Range("A24:N24").Copy()
$ws.Range("A25:N27").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Cells.Item(25,7).PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Cells.Item(26,7).PasteSpecial(-4122)
$ws.Range("G10").Copy()
$ws.Cells.Item(27,7).PasteSpecial(-4122)

$ws.Cells.Item(25,1).Value = "JAE-YONG"
$ws.Cells.Item(25,2).Value = "LILY"
$ws.Cells.Item(25,3).Value = "SANDY"
$ws.Cells.Item(25,4).Value = "TARA"
$ws.Cells.Item(25,5).Value = "SHADE"
$ws.Cells.Item(25,6).Value = "CHARLIE"
$ws.Cells.Item(25,7).Value = "Equipo 1"
$ws.Cells.Item(25,8).Value = "KRO|Nanoxx"
$ws.Cells.Item(25,9).Value = "YT : GuGu"
$ws.Cells.Item(25,10).Value = "SK|Ope"
$ws.Cells.Item(25,11).Value = "NXT|Rup"
$ws.Cells.Item(25,12).Value = "HMB|BosS"
$ws.Cells.Item(25,13).Value = "SK|Yoshi825"
$ws.Cells.Item(25,14).Value = "20250724T202642.000Z"
$ws.Cells.Item(26,1).Value = "JAE-YONG"
$ws.Cells.Item(26,2).Value = "LILY"
$ws.Cells.Item(26,3).Value = "SANDY"
$ws.Cells.Item(26,4).Value = "TARA"
$ws.Cells.Item(26,5).Value = "SHADE"
$ws.Cells.Item(26,6).Value = "CHARLIE"
$ws.Cells.Item(26,7).Value = "Equipo 1"
$ws.Cells.Item(26,8).Value = "KRO|Nanoxx"
$ws.Cells.Item(26,9).Value = "YT : GuGu"
$ws.Cells.Item(26,10).Value = "SK|Ope"
$ws.Cells.Item(26,11).Value = "NXT|Rup"
$ws.Cells.Item(26,12).Value = "HMB|BosS"
$ws.Cells.Item(26,13).Value = "SK|Yoshi825"
$ws.Cells.Item(26,14).Value = "20250724T202428.000Z"
$ws.Cells.Item(27,1).Value = "JAE-YONG"
$ws.Cells.Item(27,2).Value = "LILY"
$ws.Cells.Item(27,3).Value = "SANDY"
$ws.Cells.Item(27,4).Value = "TARA"
$ws.Cells.Item(27,5).Value = "SHADE"
$ws.Cells.Item(27,6).Value = "CHARLIE"
$ws.Cells.Item(27,7).Value = "Equipo 2"
$ws.Cells.Item(27,8).Value = "KRO|Nanoxx"
$ws.Cells.Item(27,9).Value = "YT : GuGu"
$ws.Cells.Item(27,10).Value = "SK|Ope"
$ws.Cells.Item(27,11).Value = "NXT|Rup"
$ws.Cells.Item(27,12).Value = "HMB|BosS"
$ws.Cells.Item(27,13).Value = "SK|Yoshi825"
$ws.Cells.Item(27,14).Value = "20250724T202128.000Z"

# ---- Sheet: New Horizons ----
$ws = $wb.Worksheets.Item("New Horizons")
$ws.Range("A64:N64").Copy()
$ws.Range("A65:N66").PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Cells.Item(65,7).PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Cells.Item(66,7).PasteSpecial(-4122)

$ws.Cells.Item(65,1).Value = "SPROUT"
$ws.Cells.Item(65,2).Value = "GUS"
$ws.Cells.Item(65,3).Value = "BONNIE"
$ws.Cells.Item(65,4).Value = "JANET"
$ws.Cells.Item(65,5).Value = "DOUG"
$ws.Cells.Item(65,6).Value = "R-T"
$ws.Cells.Item(65,7).Value = "Equipo 1"
$ws.Cells.Item(65,8).Value = "NOVO|Marco"
$ws.Cells.Item(65,9).Value = "TTM|Maru"
$ws.Cells.Item(65,10).Value = "TTM|Maury"
$ws.Cells.Item(65,11).Value = "Enraged 💔"
$ws.Cells.Item(65,12).Value = "SUP|Tomzy"
$ws.Cells.Item(65,13).Value = "Drage🍥"
$ws.Cells.Item(65,14).Value = "20250724T200251.000Z"
$ws.Cells.Item(66,1).Value = "SPROUT"
$ws.Cells.Item(66,2).Value = "GUS"
$ws.Cells.Item(66,3).Value = "BONNIE"
$ws.Cells.Item(66,4).Value = "JANET"
$ws.Cells.Item(66,5).Value = "DOUG"
$ws.Cells.Item(66,6).Value = "R-T"
$ws.Cells.Item(66,7).Value = "Equipo 1"
$ws.Cells.Item(66,8).Value = "NOVO|Marco"
$ws.Cells.Item(66,9).Value = "TTM|Maru"
$ws.Cells.Item(66,10).Value = "TTM|Maury"
$ws.Cells.Item(66,11).Value = "Enraged 💔"
$ws.Cells.Item(66,12).Value = "SUP|Tomzy"
$ws.Cells.Item(66,13).Value = "Drage🍥"
$ws.Cells.Item(66,14).Value = "20250724T200112.000Z"

# ---- Sheet: Open Business ----
$ws = $wb.Worksheets.Item("Open Business")
$ws.Range("A63:N63").Copy()
$ws.Range("A64:N66").PasteSpecial(-4122)
$ws.Range("G10").Copy()
$ws.Cells.Item(64,7).PasteSpecial(-4122)
$ws.Range("G4").Copy()
$ws.Cells.Item(65,7).PasteSpecial(-4122)
$ws.Range("G10").Copy()
$ws.Cells.Item(66,7).PasteSpecial(-4122)

$ws.Cells.Item(64,1).Value = "HANK"
$ws.Cells.Item(64,2).Value = "GRAY"
$ws.Cells.Item(64,3).Value = "LOU"
$ws.Cells.Item(64,4).Value = "KAZE"
$ws.Cells.Item(64,5).Value = "ASH"
$ws.Cells.Item(64,6).Value = "MOE"
$ws.Cells.Item(64,7).Value = "Equipo 1"
$ws.Cells.Item(64,8).Value = "NHG|Xemp"
$ws.Cells.Item(64,9).Value = "KCP|Fade"
$ws.Cells.Item(64,10).Value = "KCP|Tyrant"
$ws.Cells.Item(64,11).Value = "TRB|Lxffy"
$ws.Cells.Item(64,12).Value = "TRB|R B M"
$ws.Cells.Item(64,13).Value = "TRB|Zeus 解開"
$ws.Cells.Item(64,14).Value = "20250724T203010.000Z"
$ws.Cells.Item(65,1).Value = "HANK"
$ws.Cells.Item(65,2).Value = "GRAY"
$ws.Cells.Item(65,3).Value = "LOU"
$ws.Cells.Item(65,4).Value = "KAZE"
$ws.Cells.Item(65,5).Value = "ASH"
$ws.Cells.Item(65,6).Value = "MOE"
$ws.Cells.Item(65,7).Value = "Equipo 2"
$ws.Cells.Item(65,8).Value = "NHG|Xemp"
$ws.Cells.Item(65,9).Value = "KCP|Fade"
$ws.Cells.Item(65,10).Value = "KCP|Tyrant"
$ws.Cells.Item(65,11).Value = "TRB|Lxffy"
$ws.Cells.Item(65,12).Value = "TRB|R B M"
$ws.Cells.Item(65,13).Value = "TRB|Zeus 解開"
$ws.Cells.Item(65,14).Value = "20250724T202813.000Z"
$ws.Cells.Item(66,1).Value = "HANK"
$ws.Cells.Item(66,2).Value = "GRAY"
$ws.Cells.Item(66,3).Value = "LOU"
$ws.Cells.Item(66,4).Value = "KAZE"
$ws.Cells.Item(66,5).Value = "ASH"
$ws.Cells.Item(66,6).Value = "MOE"
$ws.Cells.Item(66,7).Value = "Equipo 1"
$ws.Cells.Item(66,8).Value = "NHG|Xemp"
$ws.Cells.Item(66,9).Value = "KCP|Fade"
$ws.Cells.Item(66,10).Value = "KCP|Tyrant"
$ws.Cells.Item(66,11).Value = "TRB|Lxffy"
$ws.Cells.Item(66,12).Value = "TRB|R B M"
$ws.Cells.Item(66,13).Value = "TRB|Zeus 解開"
$ws.Cells.Item(66,14).Value = "20250724T202553.000Z"
